# Line_Setup.xlsx update
#  - Row 2 (M2:T2): values that were stored as text become real numbers
#  - Row 3: new "ahpallstationpage" line-approval entry is appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to keep a text value even when it "looks" numeric,
    # then drop the temporary Text number format so no style residue remains.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Row 2: M2:T2 become numeric values instead of text ---
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 30
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 32
$ws.Range("R2").Value = 32
$ws.Range("S2").Value = 22
$ws.Range("T2").Value = 22

# --- Row 3: new entry ---
$ws.Range("A3").Value = "2025-01-29T18:57"
$ws.Range("B3").Value = "SHIFT_1"
Set-TextValue "C3" "10"
Set-TextValue "D3" "22"
Set-TextValue "E3" "33"
Set-TextValue "F3" "44"
$ws.Range("G3").Value = "OK"
$ws.Range("H3").Value = "OK"
$ws.Range("I3").Value = "OK"
$ws.Range("J3").Value = "OK"
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "OK"
Set-TextValue "M3" "1"
Set-TextValue "N3" "1"
Set-TextValue "O3" "1"
Set-TextValue "P3" "1"
Set-TextValue "Q3" "1"
Set-TextValue "R3" "1"
Set-TextValue "S3" "1"
Set-TextValue "T3" "1"
$ws.Range("U3").Value = "OK"
$ws.Range("V3").Value = "OK"
$ws.Range("W3").Value = "Test "
$ws.Range("X3").Value = "Suriya"
$ws.Range("Y3").Value = "suriya"
